$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 82 (id 80) with new match data ---
$ws.Range("B82").Value = 5574442
$ws.Range("F82").Value = "FK Qarabag"
$ws.Range("G82").Value = "FK Sumqayit"
$ws.Range("H82").Value = 1
$ws.Range("I82").Value = 2
$ws.Range("J82").Value = "A"
$ws.Range("K82").Value = 1.125
$ws.Range("L82").Value = 7.5
$ws.Range("M82").Value = 15
$ws.Range("N82").Value = 1.2
$ws.Range("O82").Value = 6
$ws.Range("P82").Value = 11
$ws.Range("Q82").Value = -2.25
$ws.Range("R82").Value = 1.975
$ws.Range("S82").Value = 1.825
$ws.Range("T82").Value = 3.5
$ws.Range("U82").Value = 1.825
$ws.Range("V82").Value = 1.975
$ws.Range("W82").Value = -1
$ws.Range("X82").Value = -1
$ws.Range("Y82").Value = 10
$ws.Range("Z82").Value = -1
$ws.Range("AA82").Value = 0.825
$ws.Range("AB82").Value = -1
$ws.Range("AC82").Value = 0.9750000000000001

# --- Update existing row 84 (id 82) with new match data ---
$ws.Range("B84").Value = 5573343
$ws.Range("F84").Value = "Shamakhi FK"
$ws.Range("G84").Value = "FK Gabala"
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = "D"
$ws.Range("K84").Value = 3.5
$ws.Range("L84").Value = 3.1
$ws.Range("M84").Value = 2
$ws.Range("N84").Value = 3.3
$ws.Range("O84").Value = 3.2
$ws.Range("P84").Value = 2.05
$ws.Range("Q84").Value = 0.25
$ws.Range("R84").Value = 2
$ws.Range("S84").Value = 1.8
$ws.Range("T84").Value = 2.5
$ws.Range("U84").Value = 1.975
$ws.Range("V84").Value = 1.825
$ws.Range("W84").Value = -1
$ws.Range("X84").Value = 2.2
$ws.Range("Y84").Value = -1
$ws.Range("Z84").Value = 0.5
$ws.Range("AA84").Value = -0.5
$ws.Range("AB84").Value = -1
$ws.Range("AC84").Value = 0.825

# --- Update existing row 85 (id 83) with new match data ---
$ws.Range("B85").Value = 5573342
$ws.Range("F85").Value = "PFK Turan Tovuz"
$ws.Range("G85").Value = "Sabail FC"
$ws.Range("H85").Value = 2
$ws.Range("I85").Value = 2
$ws.Range("J85").Value = "D"
$ws.Range("K85").Value = 2.6
$ws.Range("L85").Value = 3
$ws.Range("M85").Value = 2.6
$ws.Range("N85").Value = 2.8
$ws.Range("O85").Value = 2.875
$ws.Range("P85").Value = 2.5
$ws.Range("Q85").Value = 0
$ws.Range("R85").Value = 2.05
$ws.Range("S85").Value = 1.75
$ws.Range("T85").Value = 2.25
$ws.Range("U85").Value = 1.875
$ws.Range("V85").Value = 1.925
$ws.Range("W85").Value = -1
$ws.Range("X85").Value = 1.875
$ws.Range("Y85").Value = -1
$ws.Range("Z85").Value = 0
$ws.Range("AA85").Value = 0
$ws.Range("AB85").Value = 0.875
$ws.Range("AC85").Value = -1

# --- Update existing row 200 (id 198) with corrected/updated odds ---
$ws.Range("B200").Value = 7011605
$ws.Range("E200").Value = 45340.35416666666
$ws.Range("F200").Value = "FK Sumqayit"
$ws.Range("G200").Value = "PFK Turan Tovuz"
$ws.Range("K200").Value = 2.375
$ws.Range("L200").Value = 2.8
$ws.Range("M200").Value = 3
$ws.Range("N200").Value = 2.375
$ws.Range("O200").Value = 2.8
$ws.Range("P200").Value = 3
$ws.Range("Q200").Value = -0.25
$ws.Range("R200").Value = 2.05
$ws.Range("S200").Value = 1.75
$ws.Range("T200").Value = 2.25
$ws.Range("U200").Value = 2
$ws.Range("V200").Value = 1.8
$ws.Range("W200").Value = 0
$ws.Range("X200").Value = 0
$ws.Range("Y200").Value = 0
$ws.Range("Z200").Value = 0
$ws.Range("AA200").Value = 0

# --- Append new row 201 (id 199) for the new upcoming match ---
# Copy formatting (cell styles) from row 200 first, then fill in values,
# then clear the cells that should stay empty (result not played yet).
$ws.Range("A200:AC200").Copy()
$ws.Range("A201:AC201").PasteSpecial(-4122)

$ws.Range("A201").Value = 199
$ws.Range("B201").Value = 7011606
$ws.Range("C201").Value = "Azerbaijan Premier League"
$ws.Range("D201").Value = "Azerbaijan Premier League"
$ws.Range("E201").Value = 45340.45833333334
$ws.Range("F201").Value = "Neftchi Baku"
$ws.Range("G201").Value = "Zira IK"
$ws.Range("K201").Value = 2.2
$ws.Range("L201").Value = 2.9
$ws.Range("M201").Value = 3.2
$ws.Range("N201").Value = 2.2
$ws.Range("O201").Value = 2.9
$ws.Range("P201").Value = 3.2
$ws.Range("Q201").Value = -0.25
$ws.Range("R201").Value = 1.95
$ws.Range("S201").Value = 1.85
$ws.Range("T201").Value = 2.25
$ws.Range("U201").Value = 2
$ws.Range("V201").Value = 1.8
$ws.Range("W201").Value = 0
$ws.Range("X201").Value = 0
$ws.Range("Y201").Value = 0
$ws.Range("Z201").Value = 0
$ws.Range("AA201").Value = 0

$ws.Range("H201:J201").ClearContents()
$ws.Range("AB201:AC201").ClearContents()
